$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 149 and row 150 (columns B:G)
$tmp0 = $ws.Range("B149:G149").Value()
$ws.Range("B149:G149").Value = $ws.Range("B150:G150").Value()
$ws.Range("B150:G150").Value = $tmp0

# Swap row 264 and row 265 (columns B:G)
$tmp1 = $ws.Range("B264:G264").Value()
$ws.Range("B264:G264").Value = $ws.Range("B265:G265").Value()
$ws.Range("B265:G265").Value = $tmp1

# Swap row 279 and row 280 (columns B:G)
$tmp2 = $ws.Range("B279:G279").Value()
$ws.Range("B279:G279").Value = $ws.Range("B280:G280").Value()
$ws.Range("B280:G280").Value = $tmp2

# Swap row 313 and row 314 (columns B:G)
$tmp3 = $ws.Range("B313:G313").Value()
$ws.Range("B313:G313").Value = $ws.Range("B314:G314").Value()
$ws.Range("B314:G314").Value = $tmp3

# Swap row 316 and row 318 (columns B:G)
$tmp4 = $ws.Range("B316:G316").Value()
$ws.Range("B316:G316").Value = $ws.Range("B318:G318").Value()
$ws.Range("B318:G318").Value = $tmp4

# Swap row 346 and row 347 (columns B:G)
$tmp5 = $ws.Range("B346:G346").Value()
$ws.Range("B346:G346").Value = $ws.Range("B347:G347").Value()
$ws.Range("B347:G347").Value = $tmp5

# Swap row 372 and row 373 (columns B:G)
$tmp6 = $ws.Range("B372:G372").Value()
$ws.Range("B372:G372").Value = $ws.Range("B373:G373").Value()
$ws.Range("B373:G373").Value = $tmp6

# Swap row 379 and row 380 (columns B:G)
$tmp7 = $ws.Range("B379:G379").Value()
$ws.Range("B379:G379").Value = $ws.Range("B380:G380").Value()
$ws.Range("B380:G380").Value = $tmp7

# Swap row 382 and row 383 (columns B:G)
$tmp8 = $ws.Range("B382:G382").Value()
$ws.Range("B382:G382").Value = $ws.Range("B383:G383").Value()
$ws.Range("B383:G383").Value = $tmp8

# Swap row 389 and row 390 (columns B:G)
$tmp9 = $ws.Range("B389:G389").Value()
$ws.Range("B389:G389").Value = $ws.Range("B390:G390").Value()
$ws.Range("B390:G390").Value = $tmp9

# Swap row 400 and row 401 (columns B:G)
$tmp10 = $ws.Range("B400:G400").Value()
$ws.Range("B400:G400").Value = $ws.Range("B401:G401").Value()
$ws.Range("B401:G401").Value = $tmp10

# Swap row 419 and row 420 (columns B:G)
$tmp11 = $ws.Range("B419:G419").Value()
$ws.Range("B419:G419").Value = $ws.Range("B420:G420").Value()
$ws.Range("B420:G420").Value = $tmp11

# Swap row 421 and row 422 (columns B:G)
$tmp12 = $ws.Range("B421:G421").Value()
$ws.Range("B421:G421").Value = $ws.Range("B422:G422").Value()
$ws.Range("B422:G422").Value = $tmp12

# Swap row 431 and row 432 (columns B:G)
$tmp13 = $ws.Range("B431:G431").Value()
$ws.Range("B431:G431").Value = $ws.Range("B432:G432").Value()
$ws.Range("B432:G432").Value = $tmp13

# Swap row 536 and row 537 (columns B:G)
$tmp14 = $ws.Range("B536:G536").Value()
$ws.Range("B536:G536").Value = $ws.Range("B537:G537").Value()
$ws.Range("B537:G537").Value = $tmp14

# Swap row 579 and row 580 (columns B:G)
$tmp15 = $ws.Range("B579:G579").Value()
$ws.Range("B579:G579").Value = $ws.Range("B580:G580").Value()
$ws.Range("B580:G580").Value = $tmp15

# Swap row 583 and row 584 (columns B:G)
$tmp16 = $ws.Range("B583:G583").Value()
$ws.Range("B583:G583").Value = $ws.Range("B584:G584").Value()
$ws.Range("B584:G584").Value = $tmp16

# Swap row 590 and row 591 (columns B:G)
$tmp17 = $ws.Range("B590:G590").Value()
$ws.Range("B590:G590").Value = $ws.Range("B591:G591").Value()
$ws.Range("B591:G591").Value = $tmp17

# Swap row 604 and row 605 (columns B:G)
$tmp18 = $ws.Range("B604:G604").Value()
$ws.Range("B604:G604").Value = $ws.Range("B605:G605").Value()
$ws.Range("B605:G605").Value = $tmp18

# Swap row 687 and row 688 (columns B:G)
$tmp19 = $ws.Range("B687:G687").Value()
$ws.Range("B687:G687").Value = $ws.Range("B688:G688").Value()
$ws.Range("B688:G688").Value = $tmp19

# Swap row 720 and row 721 (columns B:G)
$tmp20 = $ws.Range("B720:G720").Value()
$ws.Range("B720:G720").Value = $ws.Range("B721:G721").Value()
$ws.Range("B721:G721").Value = $tmp20

# Swap row 859 and row 860 (columns B:F); column G is a data quirk in the
# source edit: both rows end up with the SAME new G value (13849.44), not swapped.
$tmp21 = $ws.Range("B859:F859").Value()
$ws.Range("B859:F859").Value = $ws.Range("B860:F860").Value()
$ws.Range("B860:F860").Value = $tmp21
$ws.Range("G859").Value = 13849.44
$ws.Range("G860").Value = 13849.44

# Rotate rows 350 -> 351 -> 352 -> 350 (columns B:G): new350=old351, new351=old352, new352=old350
$tmpRot = $ws.Range("B350:G350").Value()
$ws.Range("B350:G350").Value = $ws.Range("B351:G351").Value()
$ws.Range("B351:G351").Value = $ws.Range("B352:G352").Value()
$ws.Range("B352:G352").Value = $tmpRot
